$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.458.69"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.96%  '
$ws.Range('D3').Value = "'2.445.29"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.46%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'556.90"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.78%  '
$ws.Range('D6').Value = "'138.93"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.91%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'0.571"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.34%  '
$ws.Range('E9').Value = '  +3.70%  '
$ws.Range('D10').Value = "'5.80"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.02%  '
$ws.Range('D11').Value = "'0.361"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = "'24.93"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.84%  '
$ws.Range('D14').Value = "'2.878.83"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.40%  '
$ws.Range('D15').Value = "'60.402.72"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('E16').Value = '  +3.72%  '
$ws.Range('D17').Value = "'2.442.81"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').Value = "'11.43"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.19%  '
$ws.Range('D19').Value = "'4.43"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.11%  '
$ws.Range('D20').Value = "'336.12"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').Value = "'6.91"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.54%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = "'64.78"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.18%  '
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = "'1.37"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('D28').Value = "'0.0₃0794"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.12%  '
$ws.Range('E29').Value = '  +3.18%  '
$ws.Range('D30').Value = "'171.23"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('D31').Value = "'6.30"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').Value = "'18.85"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('E33').Value = '  -2.03%  '
$ws.Range('D35').Value = "'1.32"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.90%  '
$ws.Range('D36').Value = "'4.28"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = "'1.64"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = "'40.13"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('D40').Value = "'0.419"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.38%  '
$ws.Range('D41').Value = "'317.64"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.75%  '
$ws.Range('D42').Value = "'143.94"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').Value = "'3.73"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.76%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = "'0.0964"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.42%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'19.89"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.18%  '
$ws.Range('D46').Value = "'0.0525"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.32%  '
$ws.Range('D47').Value = "'0.575"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('E48').Value = '  +7.88%  '
$ws.Range('D49').Value = "'0.0227"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.12%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = "'1.64"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.00%  '
